# Updated cryptos list on Sun Nov 10 06:36:59 UTC 2024 with GitHub Actions
#
# This script refreshes the Price (column D) and Volume(1h) (column E)
# values for the crypto rows in the sheet, and fixes the ranking order of
# the "Aave" / "USDe" rows (43 and 44) which swapped places.
#
# Note: several Price values are plain numeric-looking strings (e.g.
# "205.11", "1.00", "0.999") that Excel would otherwise silently convert
# to real numbers (losing formatting, e.g. trailing zeros). Those are
# written with a leading apostrophe to force them to stay text, matching
# the source data which stores every Price cell as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "78.814.16"
$ws.Range("E2").Value = "  +2.96%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.170.01"
$ws.Range("E3").Value = "  +4.08%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - Solana
$ws.Range("D5").Value = "'205.11"
$ws.Range("E5").Value = "  +1.87%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'629.31"
$ws.Range("E6").Value = "  -0.31%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").Value = "'0.225"
$ws.Range("E8").Value = "  +10.89%  "

# Row 9
$ws.Range("D9").Value = "'0.581"
$ws.Range("E9").Value = "  +5.26%  "

# Row 10
$ws.Range("D10").Value = "3.169.67"
$ws.Range("E10").Value = "  +4.00%  "

# Row 11
$ws.Range("D11").Value = "'0.579"
$ws.Range("E11").Value = "  +32.84%  "

# Row 12
$ws.Range("D12").Value = "'0.164"
$ws.Range("E12").Value = "  +2.17%  "

# Row 13
$ws.Range("E13").Value = "  +6.40%  "

# Row 14
$ws.Range("D14").Value = "3.752.15"
$ws.Range("E14").Value = "  +3.97%  "

# Row 15
$ws.Range("D15").Value = "'0.0000223"
$ws.Range("E15").Value = "  +17.48%  "

# Row 16
$ws.Range("D16").Value = "'31.48"
$ws.Range("E16").Value = "  +7.09%  "

# Row 17
$ws.Range("D17").Value = "78.795.27"
$ws.Range("E17").Value = "  +2.96%  "

# Row 18
$ws.Range("D18").Value = "3.169.70"
$ws.Range("E18").Value = "  +4.21%  "

# Row 19
$ws.Range("D19").Value = "'14.38"
$ws.Range("E19").Value = "  +5.90%  "

# Row 20
$ws.Range("D20").Value = "'9.36"
$ws.Range("E20").Value = "  +3.50%  "

# Row 21
$ws.Range("D21").Value = "'426.88"
$ws.Range("E21").Value = "  +13.23%  "

# Row 22
$ws.Range("D22").Value = "'2.83"
$ws.Range("E22").Value = "  +24.08%  "

# Row 23
$ws.Range("D23").Value = "'4.91"
$ws.Range("E23").Value = "  +12.73%  "

# Row 24
$ws.Range("D24").Value = "'6.78"
$ws.Range("E24").Value = "  +5.07%  "

# Row 25
$ws.Range("D25").Value = "3.330.31"
$ws.Range("E25").Value = "  +3.86%  "

# Row 26
$ws.Range("D26").Value = "'4.73"
$ws.Range("E26").Value = "  +6.90%  "

# Row 27
$ws.Range("D27").Value = "'75.69"
$ws.Range("E27").Value = "  +3.28%  "

# Row 28
$ws.Range("D28").Value = "'10.91"
$ws.Range("E28").Value = "  +10.43%  "

# Row 29
$ws.Range("E29").Value = "  +0.18%  "

# Row 30
$ws.Range("E30").Value = "  +3.55%  "

# Row 31
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.21%  "

# Row 32
$ws.Range("D32").Value = "'8.85"
$ws.Range("E32").Value = "  +6.49%  "

# Row 33
$ws.Range("D33").Value = "'1.46"
$ws.Range("E33").Value = "  +4.34%  "

# Row 34
$ws.Range("D34").Value = "'510.18"
$ws.Range("E34").Value = "  -0.42%  "

# Row 35
$ws.Range("E35").Value = "  +0.37%  "

# Row 36
$ws.Range("D36").Value = "'0.127"
$ws.Range("E36").Value = "  +21.13%  "

# Row 37
$ws.Range("D37").Value = "'22.85"
$ws.Range("E37").Value = "  +9.39%  "

# Row 38
$ws.Range("E38").Value = "  +19.40%  "

# Row 39
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  -0.03%  "

# Row 40
$ws.Range("D40").Value = "'0.396"
$ws.Range("E40").Value = "  +3.23%  "

# Row 41
$ws.Range("D41").Value = "'163.37"
$ws.Range("E41").Value = "  -0.56%  "

# Row 42
$ws.Range("D42").Value = "'19.96"
$ws.Range("E42").Value = "  -0.27%  "

# Row 43 - was Aave, now USDe (ranking order changed)
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.13%  "

# Row 44 - was USDe, now Aave (ranking order changed)
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'191.93"
$ws.Range("E44").Value = "  +0.34%  "

# Row 45
$ws.Range("D45").Value = "'5.37"
$ws.Range("E45").Value = "  +5.69%  "

# Row 46
$ws.Range("D46").Value = "'0.807"
$ws.Range("E46").Value = "  +14.13%  "

# Row 47
$ws.Range("D47").Value = "'1.78"
$ws.Range("E47").Value = "  +7.14%  "

# Row 48
$ws.Range("E48").Value = "  +2.73%  "

# Row 49
$ws.Range("D49").Value = "'42.52"
$ws.Range("E49").Value = "  -1.45%  "

# Row 50
$ws.Range("D50").Value = "'2.49"
$ws.Range("E50").Value = "  +5.93%  "

# Row 51
$ws.Range("E51").Value = "  +1.91%  "
